$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows, continuing the series after "01-10-2021" (row 275).
# Each new row uses the same B/C/D values as the last existing row (275).
$dates = @("02-10-2021", "03-10-2021", "04-10-2021", "05-10-2021", "06-10-2021")

$startRow = 276
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i

    # Leading apostrophe forces text entry (otherwise Excel auto-parses
    # "dd-mm-yyyy" looking strings into date serials); resetting the style
    # to "Normal" afterwards drops the quote-prefix formatting flag so the
    # cell ends up with no explicit style, matching the rest of the column.
    $ws.Cells.Item($row, 1).Value = "'" + $dates[$i]
    $ws.Cells.Item($row, 1).Style = "Normal"

    $ws.Cells.Item($row, 2).Value = 12836
    $ws.Cells.Item($row, 3).Value = 266
    $ws.Cells.Item($row, 4).Value = 393
}
